$d = $word.ActiveDocument

# --- Change 1 ---
# "Uma quiz que conterá 8 questões referente ao meu site, onde ao final
# aparecera a quantidade de acertos e porcentagem, além de um quiz de
# personalidade para identificar qual personagem você se identifica"
# becomes
# "Uma quiz que conterá 5 questões referente ao meu site, onde ao final
# aparecera a quantidade de acertos e porcentagem."
$old1 = "8 questões referente ao meu site, onde ao final aparecera a quantidade de acertos e porcentagem, além de um quiz de personalidade para identificar qual personagem você se identifica"
$new1 = "5 questões referente ao meu site, onde ao final aparecera a quantidade de acertos e porcentagem."
$r1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- Change 2 ---
# "Gráfico com personagens favoritos de cada usuário e média de pontuação
# de cada usuário" becomes
# "Gráfico com personagens favoritos de cada usuário e da média de
# pontuação de cada usuário "
$old2 = "Gráfico com personagens favoritos de cada usuário e média de pontuação de cada usuário"
$new2 = "Gráfico com personagens favoritos de cada usuário e da média de pontuação de cada usuário "
$r2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# --- Change 3a ---
# "Painel com resultado do Quis" -> "Painel com resultado do Quiz"
$old3 = "Quis"
$new3 = "Quiz"
$r3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# --- Change 3b ---
# Delete the whole paragraph "Módulo de Quiz de personalidade"
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "Módulo de Quiz de personalidade") {
        $p.Range.Delete()
        break
    }
}

# --- Change 4 ---
# Add a new paragraph "Fazer App Mobile." right after the paragraph
# "Músicas temas da abertura do anime" (same list numbering, numId=4)
$fr = $d.Content
$found4 = $fr.Find.Execute("úsicas temas da abertura do anime", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $fr.Collapse(0)
    $fr.InsertAfter("`rFazer App Mobile.")
}

Write-Output "r1=$r1 r2=$r2 r3=$r3 found4=$found4"
